# Weekly update: insert a new price record as row 36 (pushing existing
# rows 36-98 down to 37-99), matching the "Fruta / hortaliza, semanal"
# commit that adds the latest week's observation to the top of the
# historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 36; everything below shifts down.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(36, 1).Value  = 1
$ws.Cells.Item(36, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(36, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(36, 4).Value  = 44797
$ws.Cells.Item(36, 5).Value  = 15
$ws.Cells.Item(36, 6).Value  = "Fruta"
$ws.Cells.Item(36, 7).Value  = 100102
$ws.Cells.Item(36, 8).Value  = "Cítricos"
$ws.Cells.Item(36, 9).Value  = 100102005
$ws.Cells.Item(36, 10).Value = "Naranja"
$ws.Cells.Item(36, 11).Value = "Fukumoto"
$ws.Cells.Item(36, 12).Value = "Segunda"
$ws.Cells.Item(36, 13).Value = 250
$ws.Cells.Item(36, 14).Value = 650
$ws.Cells.Item(36, 15).Value = 700
$ws.Cells.Item(36, 16).Value = 675
$ws.Cells.Item(36, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(36, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(36, 19).Value = 675
$ws.Cells.Item(36, 20).Value = 1
